$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

$ws.Range("H113").Value = 2487.5
$ws.Range("I113").Value = 2487.5
$ws.Range("K113").Value = 2487.5
$ws.Range("M113").Value = 766.5

$ws.Range("H132").Value = 13008.125
$ws.Range("I132").Value = 12240.25
$ws.Range("K132").Value = 36720.75
$ws.Range("M132").Value = -34190.75

$ws.Range("H137").Value = 2953.652
$ws.Range("J137").Value = 3499.2144
$ws.Range("L137").Value = 10497.6432
$ws.Range("N137").Value = -15597.6432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2704970.2
$ws.Range("I32").Value = 480.03226
$ws.Range("K32").Value = 480.03226
$ws.Range("M32").Value = -193.03226

$ws.Range("H96").Value = 3363021.5
$ws.Range("J96").Value = 3363021.5
$ws.Range("L96").Value = 3363021.5
$ws.Range("N96").Value = -3368513.5

$ws.Range("H97").Value = 1761.4166
$ws.Range("I97").Value = 1679.375
$ws.Range("K97").Value = 1679.375
$ws.Range("M97").Value = -1183.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 9882927
$ws.Range("I7").Value = 11500180
$ws.Range("K7").Value = 11500180
$ws.Range("M7").Value = -11500067

$ws.Range("H94").Value = 420
$ws.Range("I94").Value = 420
$ws.Range("K94").Value = 420
$ws.Range("M94").Value = 31

$ws.Range("H99").Value = 76923940
$ws.Range("I99").Value = 90909910
$ws.Range("K99").Value = 90909910
$ws.Range("M99").Value = -90908412

$ws.Range("H134").Value = 2653.9487
$ws.Range("I134").Value = 880.4838999999999
$ws.Range("K134").Value = 2641.4517
$ws.Range("M134").Value = -106.4516999999996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5005.3096
$ws.Range("I31").Value = 2068.5833
$ws.Range("J31").Value = 6180
$ws.Range("K31").Value = 2068.5833
$ws.Range("L31").Value = 6180
$ws.Range("M31").Value = -1773.5833
$ws.Range("N31").Value = -6770

$ws.Range("H34").Value = 5005.3096
$ws.Range("I34").Value = 2068.5833
$ws.Range("J34").Value = 6180
$ws.Range("K34").Value = 2068.5833
$ws.Range("L34").Value = 6180
$ws.Range("M34").Value = -1866.5833
$ws.Range("N34").Value = -6584

$ws.Range("H107").Value = 857.1111
$ws.Range("I107").Value = 239.66667
$ws.Range("K107").Value = 239.66667
$ws.Range("M107").Value = 1680.33333

$ws.Range("H125").Value = 70575
$ws.Range("J125").Value = 70575
$ws.Range("L125").Value = 70575
$ws.Range("N125").Value = -75495

$ws.Range("H138").Value = 126996.8
$ws.Range("J138").Value = 126996.8
$ws.Range("L138").Value = 126996.8
$ws.Range("N138").Value = -137276.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 23.457144
$ws.Range("I2").Value = 23.8
$ws.Range("J2").Value = 23
$ws.Range("K2").Value = 142.8
$ws.Range("L2").Value = 138
$ws.Range("M2").Value = -29.80000000000001
$ws.Range("N2").Value = -364

$ws.Range("H17").Value = 433.33334
$ws.Range("J17").Value = 800
$ws.Range("L17").Value = 2400
$ws.Range("N17").Value = -2738

$ws.Range("H39").Value = 6578.273
$ws.Range("J39").Value = 6578.273
$ws.Range("L39").Value = 19734.819
$ws.Range("N39").Value = -20322.819

$ws.Range("H55").Value = 3477.5186
$ws.Range("J55").Value = 3592.423
$ws.Range("L55").Value = 10777.269
$ws.Range("N55").Value = -11131.269

$ws.Range("H68").Value = 694
$ws.Range("J68").Value = 701.5
$ws.Range("L68").Value = 2104.5
$ws.Range("N68").Value = -3726.5

$ws.Range("H71").Value = 694
$ws.Range("J71").Value = 701.5
$ws.Range("L71").Value = 6313.5
$ws.Range("N71").Value = -14425.5

$ws.Range("H113").Value = 1830.9166
$ws.Range("I113").Value = 1050
$ws.Range("J113").Value = 1987.1
$ws.Range("K113").Value = 3150
$ws.Range("L113").Value = 5961.299999999999
$ws.Range("M113").Value = -980
$ws.Range("N113").Value = -10301.3

$ws.Range("H140").Value = 3896.3333
$ws.Range("I140").Value = 3298.75
$ws.Range("K140").Value = 9896.25
$ws.Range("M140").Value = -4716.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 238.33333
$ws.Range("I9").Value = 238.33333
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 238.33333
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -68.33332999999999
$ws.Range("N9").ClearContents()

$ws.Range("H80").Value = 4187.25
$ws.Range("I80").Value = 4329.6
$ws.Range("J80").Value = 3950
$ws.Range("K80").Value = 4329.6
$ws.Range("L80").Value = 3950
$ws.Range("M80").Value = -3331.6
$ws.Range("N80").Value = -5946

$ws.Range("H83").Value = 4187.25
$ws.Range("I83").Value = 4329.6
$ws.Range("J83").Value = 3950
$ws.Range("K83").Value = 21648
$ws.Range("L83").Value = 19750
$ws.Range("M83").Value = -16656
$ws.Range("N83").Value = -29734

$ws.Range("H113").Value = 7443.7
$ws.Range("I113").Value = 4664.3335
$ws.Range("J113").Value = 8634.857
$ws.Range("K113").Value = 4664.3335
$ws.Range("L113").Value = 8634.857
$ws.Range("M113").Value = -2494.3335
$ws.Range("N113").Value = -12974.857

$ws.Range("H126").Value = 5699.75
$ws.Range("I126").Value = 3666.3333
$ws.Range("J126").Value = 11800
$ws.Range("K126").Value = 10998.9999
$ws.Range("L126").Value = 35400
$ws.Range("M126").Value = -8528.999899999999
$ws.Range("N126").Value = -40340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 591.8570999999999
$ws.Range("I9").Value = 607.3333
$ws.Range("K9").Value = 607.3333
$ws.Range("M9").Value = -383.3333

$ws.Range("H10").Value = 1717.5
$ws.Range("J10").Value = 2575
$ws.Range("L10").Value = 2575
$ws.Range("N10").Value = -2855

$ws.Range("H68").Value = 8222.223
$ws.Range("J68").Value = 9500
$ws.Range("L68").Value = 9500
$ws.Range("N68").Value = -10998

$ws.Range("H71").Value = 8222.223
$ws.Range("J71").Value = 9500
$ws.Range("L71").Value = 47500
$ws.Range("N71").Value = -54988

$ws.Range("H132").Value = 3130.6667
$ws.Range("I132").Value = 1956.4667
$ws.Range("K132").Value = 5869.4001
$ws.Range("M132").Value = -3339.4001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7650
$ws.Range("I62").Value = 3701
$ws.Range("J62").Value = 8527.556
$ws.Range("K62").Value = 3701
$ws.Range("L62").Value = 8527.556
$ws.Range("M62").Value = -3077
$ws.Range("N62").Value = -9775.556

$ws.Range("H65").Value = 7650
$ws.Range("I65").Value = 3701
$ws.Range("J65").Value = 8527.556
$ws.Range("K65").Value = 18505
$ws.Range("L65").Value = 42637.78
$ws.Range("M65").Value = -15385
$ws.Range("N65").Value = -48877.78

$ws.Range("H141").Value = 140356.5
$ws.Range("J141").Value = 140356.5
$ws.Range("L141").Value = 140356.5
$ws.Range("N141").Value = -150716.5
